$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 'Tokiahuru at Karioi Domain Road'
$ws.Range("B2").Value = 'Chlorophyll A'
$ws.Range("C2").Value = 5
$ws.Range("D2").Value = $false
$ws.Range("E2").Value = 'ok'
$ws.Range("F2").Value = 0.989097769253366
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0.857142857142857
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 19.5
$ws.Range("K2").Value = -5.14034716109155
$ws.Range("L2").Value = -10.3221798846321
$ws.Range("M2").Value = -1.36570278324329
$ws.Range("N2").Value = -26.3607546722644
$ws.Range("O2").Value = 'RepSite'
$ws.Range("P2").Value = 'Extremely likely improving'
$ws.Range("Q2").Value = 1815033.57
$ws.Range("R2").Value = 5627502.8
$ws.Range("S2").Value = 'Ruapehu District'
$ws.Range("T2").Value = 'Whangaehu'
$ws.Range("U2").Value = 'Upper Whangaehu'
$ws.Range("V2").Value = 'Whau_1c'
$ws.Range("W2").Value = 'mg/m2'

# Row 3
$ws.Range("A3").Value = 'Tokiahuru at Karioi Domain Road'
$ws.Range("B3").Value = 'Chlorophyll A'
$ws.Range("C3").Value = 10
$ws.Range("D3").Value = $false
$ws.Range("E3").Value = 'ok'
$ws.Range("F3").Value = 0.850488543396862
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0.726027397260274
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 16.5
$ws.Range("K3").Value = -0.6670071552628
$ws.Range("L3").Value = -1.69719213330578
$ws.Range("M3").Value = 0.472898553130199
$ws.Range("N3").Value = -4.04246760765333
$ws.Range("O3").Value = 'RepSite'
$ws.Range("P3").Value = 'Likely improving'
$ws.Range("Q3").Value = 1815033.57
$ws.Range("R3").Value = 5627502.8
$ws.Range("S3").Value = 'Ruapehu District'
$ws.Range("T3").Value = 'Whangaehu'
$ws.Range("U3").Value = 'Upper Whangaehu'
$ws.Range("V3").Value = 'Whau_1c'
$ws.Range("W3").Value = 'mg/m2'

# Row 4
$ws.Range("A4").Value = 'Tokiahuru at Karioi Domain Road'
$ws.Range("B4").Value = 'Chlorophyll A'
$ws.Range("C4").Value = 15
$ws.Range("D4").Value = $false
$ws.Range("E4").Value = 'ok'
$ws.Range("F4").Value = 0.0000308299164009869
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0.645669291338583
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 11.5
$ws.Range("K4").Value = 1.02646370023419
$ws.Range("L4").Value = 0.559612289746005
$ws.Range("M4").Value = 1.54746253907928
$ws.Range("N4").Value = 8.92577130638428
$ws.Range("O4").Value = 'RepSite'
$ws.Range("P4").Value = 'Exceptionally unlikely improving'
$ws.Range("Q4").Value = 1815033.57
$ws.Range("R4").Value = 5627502.8
$ws.Range("S4").Value = 'Ruapehu District'
$ws.Range("T4").Value = 'Whangaehu'
$ws.Range("U4").Value = 'Upper Whangaehu'
$ws.Range("V4").Value = 'Whau_1c'
$ws.Range("W4").Value = 'mg/m2'

# Row 5
$ws.Range("A5").Value = 'Tokiahuru at Karioi Domain Road'
$ws.Range("B5").Value = 'ASPM (Macroinvertebrate Average Score Per Metric)'
$ws.Range("C5").Value = 5
$ws.Range("D5").Value = $false
$ws.Range("E5").Value = 'ok'
$ws.Range("F5").Value = 0.768783636774762
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 1
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0.37
$ws.Range("K5").Value = 0.0133626990772081
$ws.Range("L5").Value = -0.0663349081991328
$ws.Range("M5").Value = 0.0507271390824439
$ws.Range("N5").Value = 3.61154029113733
$ws.Range("O5").Value = 'RepSite'
$ws.Range("P5").Value = 'Likely improving'
$ws.Range("Q5").Value = 1815033.57
$ws.Range("R5").Value = 5627502.8
$ws.Range("S5").Value = 'Ruapehu District'
$ws.Range("T5").Value = 'Whangaehu'
$ws.Range("U5").Value = 'Upper Whangaehu'
$ws.Range("V5").Value = 'Whau_1c'

# Row 6
$ws.Range("A6").Value = 'Tokiahuru at Karioi Domain Road'
$ws.Range("B6").Value = 'MCI (Macroinvertebrate Community Index)'
$ws.Range("C6").Value = 5
$ws.Range("D6").Value = $false
$ws.Range("E6").Value = 'ok'
$ws.Range("F6").Value = 0.5
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 1
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 111
$ws.Range("K6").Value = -0.617317704164059
$ws.Range("L6").Value = -13.5112479791629
$ws.Range("M6").Value = 9.85063491979848
$ws.Range("N6").Value = -0.556142075823477
$ws.Range("O6").Value = 'RepSite'
$ws.Range("P6").Value = 'As likely as not improving'
$ws.Range("Q6").Value = 1815033.57
$ws.Range("R6").Value = 5627502.8
$ws.Range("S6").Value = 'Ruapehu District'
$ws.Range("T6").Value = 'Whangaehu'
$ws.Range("U6").Value = 'Upper Whangaehu'
$ws.Range("V6").Value = 'Whau_1c'

# Row 7
$ws.Range("A7").Value = 'Tokiahuru at Karioi Domain Road'
$ws.Range("B7").Value = 'QMCI (Quantitative Macroinvertebrate Community Index)'
$ws.Range("C7").Value = 5
$ws.Range("D7").Value = $false
$ws.Range("E7").Value = 'ok'
$ws.Range("F7").Value = 0.768783636774762
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 1
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 4
$ws.Range("K7").Value = 0.467889583458621
$ws.Range("L7").Value = -0.670057212467726
$ws.Range("M7").Value = 0.704927105436432
$ws.Range("N7").Value = 11.6972395864655
$ws.Range("O7").Value = 'RepSite'
$ws.Range("P7").Value = 'Likely improving'
$ws.Range("Q7").Value = 1815033.57
$ws.Range("R7").Value = 5627502.8
$ws.Range("S7").Value = 'Ruapehu District'
$ws.Range("T7").Value = 'Whangaehu'
$ws.Range("U7").Value = 'Upper Whangaehu'
$ws.Range("V7").Value = 'Whau_1c'

# Row 8
$ws.Range("A8").Value = 'Tokiahuru at Karioi Domain Road'
$ws.Range("B8").Value = 'ASPM (Macroinvertebrate Average Score Per Metric)'
$ws.Range("C8").Value = 10
$ws.Range("D8").Value = $false
$ws.Range("E8").Value = 'ok'
$ws.Range("F8").Value = 0.429013828493761
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 1
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0.3795
$ws.Range("K8").Value = -0.0077042186001917
$ws.Range("L8").Value = -0.0431795551072911
$ws.Range("M8").Value = 0.0104420062953873
$ws.Range("N8").Value = -2.03009712785026
$ws.Range("O8").Value = 'RepSite'
$ws.Range("P8").Value = 'As likely as not improving'
$ws.Range("Q8").Value = 1815033.57
$ws.Range("R8").Value = 5627502.8
$ws.Range("S8").Value = 'Ruapehu District'
$ws.Range("T8").Value = 'Whangaehu'
$ws.Range("U8").Value = 'Upper Whangaehu'
$ws.Range("V8").Value = 'Whau_1c'

# Row 9
$ws.Range("A9").Value = 'Tokiahuru at Karioi Domain Road'
$ws.Range("B9").Value = 'MCI (Macroinvertebrate Community Index)'
$ws.Range("C9").Value = 10
$ws.Range("D9").Value = $false
$ws.Range("E9").Value = 'ok'
$ws.Range("F9").Value = 0.295752518497458
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 1
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 111.5
$ws.Range("K9").Value = -1.04956896551724
$ws.Range("L9").Value = -4.0902001935726
$ws.Range("M9").Value = 0.729196334420924
$ws.Range("N9").Value = -0.941317457862997
$ws.Range("O9").Value = 'RepSite'
$ws.Range("P9").Value = 'Unlikely improving'
$ws.Range("Q9").Value = 1815033.57
$ws.Range("R9").Value = 5627502.8
$ws.Range("S9").Value = 'Ruapehu District'
$ws.Range("T9").Value = 'Whangaehu'
$ws.Range("U9").Value = 'Upper Whangaehu'
$ws.Range("V9").Value = 'Whau_1c'

# Row 10
$ws.Range("A10").Value = 'Tokiahuru at Karioi Domain Road'
$ws.Range("B10").Value = 'QMCI (Quantitative Macroinvertebrate Community Index)'
$ws.Range("C10").Value = 10
$ws.Range("D10").Value = $false
$ws.Range("E10").Value = 'ok'
$ws.Range("F10").Value = 0.429013828493761
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 1
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 3.985
$ws.Range("K10").Value = -0.0576510903426792
$ws.Range("L10").Value = -0.328637058482785
$ws.Range("M10").Value = 0.190223005577773
$ws.Range("N10").Value = -1.446702392539
$ws.Range("O10").Value = 'RepSite'
$ws.Range("P10").Value = 'As likely as not improving'
$ws.Range("Q10").Value = 1815033.57
$ws.Range("R10").Value = 5627502.8
$ws.Range("S10").Value = 'Ruapehu District'
$ws.Range("T10").Value = 'Whangaehu'
$ws.Range("U10").Value = 'Upper Whangaehu'
$ws.Range("V10").Value = 'Whau_1c'

# Row 11
$ws.Range("A11").Value = 'Tokiahuru at Karioi Domain Road'
$ws.Range("B11").Value = 'ASPM (Macroinvertebrate Average Score Per Metric)'
$ws.Range("C11").Value = 15
$ws.Range("D11").Value = $false
$ws.Range("E11").Value = 'ok'
$ws.Range("F11").Value = 0.094443738769769
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 1
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 0.4045
$ws.Range("K11").Value = -0.0173033974011911
$ws.Range("L11").Value = -0.026560404499094
$ws.Range("M11").Value = 0.0011372236777692
$ws.Range("N11").Value = -4.27772494467024
$ws.Range("O11").Value = 'RepSite'
$ws.Range("P11").Value = 'Very unlikely improving'
$ws.Range("Q11").Value = 1815033.57
$ws.Range("R11").Value = 5627502.8
$ws.Range("S11").Value = 'Ruapehu District'
$ws.Range("T11").Value = 'Whangaehu'
$ws.Range("U11").Value = 'Upper Whangaehu'
$ws.Range("V11").Value = 'Whau_1c'

# Row 12
$ws.Range("A12").Value = 'Tokiahuru at Karioi Domain Road'
$ws.Range("B12").Value = 'MCI (Macroinvertebrate Community Index)'
$ws.Range("C12").Value = 15
$ws.Range("D12").Value = $false
$ws.Range("E12").Value = 'ok'
$ws.Range("F12").Value = 0.002050700430952
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 1
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 114.4
$ws.Range("K12").Value = -2.58325645756458
$ws.Range("L12").Value = -3.54821195634238
$ws.Range("M12").Value = -1.40035936470761
$ws.Range("N12").Value = -2.25809130905994
$ws.Range("O12").Value = 'RepSite'
$ws.Range("P12").Value = 'Exceptionally unlikely improving'
$ws.Range("Q12").Value = 1815033.57
$ws.Range("R12").Value = 5627502.8
$ws.Range("S12").Value = 'Ruapehu District'
$ws.Range("T12").Value = 'Whangaehu'
$ws.Range("U12").Value = 'Upper Whangaehu'
$ws.Range("V12").Value = 'Whau_1c'

# Row 13
$ws.Range("A13").Value = 'Tokiahuru at Karioi Domain Road'
$ws.Range("B13").Value = 'QMCI (Quantitative Macroinvertebrate Community Index)'
$ws.Range("C13").Value = 15
$ws.Range("D13").Value = $false
$ws.Range("E13").Value = 'ok'
$ws.Range("F13").Value = 0.024372460369853
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 1
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 4.398
$ws.Range("K13").Value = -0.234446510036496
$ws.Range("L13").Value = -0.368577877420468
$ws.Range("M13").Value = -0.0722691156182014
$ws.Range("N13").Value = -5.33075284303084
$ws.Range("O13").Value = 'RepSite'
$ws.Range("P13").Value = 'Extremely unlikely improving'
$ws.Range("Q13").Value = 1815033.57
$ws.Range("R13").Value = 5627502.8
$ws.Range("S13").Value = 'Ruapehu District'
$ws.Range("T13").Value = 'Whangaehu'
$ws.Range("U13").Value = 'Upper Whangaehu'
$ws.Range("V13").Value = 'Whau_1c'

Write-Host "done"